$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reduce the shared-string table ("c"/"d" removed) and mirror A1:A2 into B1:B2,
# shifting the old C/D rows so A3/A4 now read "a"/"b" and A5 stays "e".
$ws.Range("B1").Value = $ws.Range("A1").Value2
$ws.Range("B2").Value = $ws.Range("A2").Value2
$ws.Range("A3").Value = $ws.Range("A1").Value2
$ws.Range("A4").Value = $ws.Range("A2").Value2

$ws.Range("B4").Select()
